$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.53"
$ws.Range("E2").Value = "'1.77%"
$ws.Range("D3").Value = "'31.83"
$ws.Range("E3").Value = "'0.47%"
$ws.Range("D4").Value = "'5.020"
$ws.Range("E4").Value = "'-0.82%"
$ws.Range("D5").Value = "'0.07808"
$ws.Range("E5").Value = "'-3.01%"
$ws.Range("D6").Value = "'2.044"
$ws.Range("E6").Value = "'-21.39%"
$ws.Range("D7").Value = "'7.781"
$ws.Range("E7").Value = "'-0.16%"
$ws.Range("D8").Value = "'3.781"
$ws.Range("E8").Value = "'-0.99%"
$ws.Range("D9").Value = "'0.9196"
$ws.Range("E9").Value = "'-0.38%"
$ws.Range("D10").Value = "'0.1744"
$ws.Range("E10").Value = "'-0.66%"
$ws.Range("D11").Value = "'0.07851"
$ws.Range("E11").Value = "'5.82%"
$ws.Range("D12").Value = "'0.08835"
$ws.Range("E12").Value = "'-0.46%"
$ws.Range("D13").Value = "'0.03127"
$ws.Range("E13").Value = "'2.96%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("D15").Value = "'0.001518"
$ws.Range("E15").Value = "'1.68%"
$ws.Range("D16").Value = "'0.005814"
$ws.Range("E16").Value = "'-4.05%"
$ws.Range("D17").Value = "'3.464"
$ws.Range("E17").Value = "'-2.47%"
$ws.Range("D18").Value = "'2.268"
$ws.Range("E18").Value = "'0.89%"
$ws.Range("E19").Value = "'1.07%"
$ws.Range("D20").Value = "'0.1292"
$ws.Range("E20").Value = "'-3.30%"
$ws.Range("D21").Value = "'4.184"
$ws.Range("E21").Value = "'4.75%"
$ws.Range("D22").Value = "'0.1809"
$ws.Range("E22").Value = "'9.75%"
$ws.Range("D23").Value = "'0.04609"
$ws.Range("E23").Value = "'0.49%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'0.13%"
$ws.Range("D25").Value = "'0.004473"
$ws.Range("E25").Value = "'0.60%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'4.44%"
$ws.Range("D39").Value = "'0.01748"
$ws.Range("E39").Value = "'-1.44%"
$ws.Range("D40").Value = "'0.04756"
$ws.Range("E40").Value = "'5.87%"
$ws.Range("D41").Value = "'0.007118"
$ws.Range("E41").Value = "'5.84%"
$ws.Range("D42").Value = "'0.1354"
$ws.Range("E42").Value = "'0.45%"
$ws.Range("D43").Value = "'0.002097"
$ws.Range("E43").Value = "'-4.95%"
$ws.Range("D44").Value = "'0.01077"
$ws.Range("E44").Value = "'9.42%"
$ws.Range("D45").Value = "'0.00006070"
$ws.Range("E45").Value = "'-6.17%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.24%"
$ws.Range("D47").Value = "'0.003552"
$ws.Range("E47").Value = "'-59.35%"
$ws.Range("D48").Value = "'1.175"
$ws.Range("E48").Value = "'43.23%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.24%"
$ws.Range("E50").Value = "'0.24%"
